$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2591.6365
$ws.Range("I100").Value = 1400
$ws.Range("J100").Value = 2710.8
$ws.Range("K100").Value = 1400
$ws.Range("L100").Value = 2710.8
$ws.Range("M100").Value = -859
$ws.Range("N100").Value = -3792.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I137").Value = 820.75
$ws.Range("J137").Value = 200000000
$ws.Range("K137").Value = 2462.25
$ws.Range("L137").Value = 600000000
$ws.Range("M137").Value = 87.75
$ws.Range("N137").Value = -600005100

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H140").Value = 48733.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 48733.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 48733.332
$ws.Range("N140").Value = -59093.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1079.7646
$ws.Range("I2").Value = 1027.9231
$ws.Range("J2").Value = 1248.25
$ws.Range("K2").Value = 1027.9231
$ws.Range("L2").Value = 1248.25
$ws.Range("M2").Value = -914.9231
$ws.Range("N2").Value = -1474.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9556.056
$ws.Range("I32").Value = 9964.121999999999
$ws.Range("J32").Value = 8269.076999999999
$ws.Range("K32").Value = 9964.121999999999
$ws.Range("L32").Value = 8269.076999999999
$ws.Range("M32").Value = -9677.121999999999
$ws.Range("N32").Value = -8843.076999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2224.4614
$ws.Range("I45").Value = 2757.7144
$ws.Range("J45").Value = 1602.3334
$ws.Range("K45").Value = 2757.7144
$ws.Range("L45").Value = 1602.3334
$ws.Range("M45").Value = -2380.7144
$ws.Range("N45").Value = -2356.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 31253142
$ws.Range("I61").Value = 31253142
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 31253142
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -31252930

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1079.7646
$ws.Range("I116").Value = 1027.9231
$ws.Range("J116").Value = 1248.25
$ws.Range("K116").Value = 1027.9231
$ws.Range("L116").Value = 1248.25
$ws.Range("M116").Value = 1266.0769
$ws.Range("N116").Value = -5836.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 31253142
$ws.Range("I136").Value = 31253142
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 93759426
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -93756876

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1079.7646
$ws.Range("I3").Value = 1027.9231
$ws.Range("J3").Value = 1248.25
$ws.Range("K3").Value = 1027.9231
$ws.Range("L3").Value = 1248.25
$ws.Range("M3").Value = -913.9231
$ws.Range("N3").Value = -1476.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1179.6923
$ws.Range("I107").Value = 1140.375
$ws.Range("J107").Value = 1242.6
$ws.Range("K107").Value = 1140.375
$ws.Range("L107").Value = 1242.6
$ws.Range("M107").Value = 779.625
$ws.Range("N107").Value = -5082.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3126.6667
$ws.Range("I62").Value = 2190
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2190
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1566
$ws.Range("N62").Value = -6248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3126.6667
$ws.Range("I65").Value = 2190
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 10950
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -7830
$ws.Range("N65").Value = -31240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 995.5714
$ws.Range("I107").Value = 745.4666999999999
$ws.Range("J107").Value = 1620.8334
$ws.Range("K107").Value = 745.4666999999999
$ws.Range("L107").Value = 1620.8334
$ws.Range("M107").Value = 1174.5333
$ws.Range("N107").Value = -5460.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 16668957
$ws.Range("I132").Value = 25001790
$ws.Range("J132").Value = 3292.3
$ws.Range("K132").Value = 75005370
$ws.Range("L132").Value = 9876.900000000001
$ws.Range("M132").Value = -75002840
$ws.Range("N132").Value = -14936.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 42489.918
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 42489.918
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 42489.918
$ws.Range("N140").Value = -52849.918

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 8744.429
$ws.Range("I63").Value = 4570.3335
$ws.Range("J63").Value = 11875
$ws.Range("K63").Value = 13711.0005
$ws.Range("L63").Value = 35625
$ws.Range("M63").Value = -12962.0005
$ws.Range("N63").Value = -37123

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H66").Value = 8744.429
$ws.Range("I66").Value = 4570.3335
$ws.Range("J66").Value = 11875
$ws.Range("K66").Value = 41133.0015
$ws.Range("L66").Value = 106875
$ws.Range("M66").Value = -37389.0015
$ws.Range("N66").Value = -114363

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4556.5
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 5136
$ws.Range("K70").Value = 1500
$ws.Range("L70").Value = 15408
$ws.Range("M70").Value = -1185
$ws.Range("N70").Value = -16038

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 4556.5
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 5136
$ws.Range("K73").Value = 1500
$ws.Range("L73").Value = 15408
$ws.Range("M73").Value = -408
$ws.Range("N73").Value = -17592

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1617.091
$ws.Range("I118").Value = 1600
$ws.Range("J118").Value = 1620.8889
$ws.Range("K118").Value = 4800
$ws.Range("L118").Value = 4862.6667
$ws.Range("M118").Value = -3557
$ws.Range("N118").Value = -7348.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 834.7
$ws.Range("I131").Value = 438.33334
$ws.Range("J131").Value = 860
$ws.Range("K131").Value = 1315.00002
$ws.Range("L131").Value = 2580
$ws.Range("M131").Value = 3724.99998
$ws.Range("N131").Value = -12660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4446807.5
$ws.Range("I122").Value = 9525412
$ws.Range("J122").Value = 3028.25
$ws.Range("K122").Value = 28576236
$ws.Range("L122").Value = 9084.75
$ws.Range("M122").Value = -28573786
$ws.Range("N122").Value = -13984.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 40319.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 40319.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 40319.5
$ws.Range("N123").Value = -45219.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 889.375
$ws.Range("I22").Value = 605.55554
$ws.Range("J22").Value = 1254.2858
$ws.Range("K22").Value = 605.55554
$ws.Range("L22").Value = 1254.2858
$ws.Range("M22").Value = -310.55554
$ws.Range("N22").Value = -1844.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 889.375
$ws.Range("I27").Value = 605.55554
$ws.Range("J27").Value = 1254.2858
$ws.Range("K27").Value = 605.55554
$ws.Range("L27").Value = 1254.2858
$ws.Range("M27").Value = -498.55554
$ws.Range("N27").Value = -1468.2858

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I40").Value = 6390.0835
$ws.Range("J40").Value = 6233.3335
$ws.Range("K40").Value = 6390.0835
$ws.Range("L40").Value = 6233.3335
$ws.Range("M40").Value = -6254.0835
$ws.Range("N40").Value = -6505.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 37653
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 37653
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 37653
$ws.Range("N81").Value = -39649

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H84").Value = 37653
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 37653
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 112959
$ws.Range("N84").Value = -122943

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9556.817999999999
$ws.Range("I132").Value = 6949.385
$ws.Range("J132").Value = 11251.65
$ws.Range("K132").Value = 20848.155
$ws.Range("L132").Value = 33754.95
$ws.Range("M132").Value = -18318.155
$ws.Range("N132").Value = -38814.95

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 59963.43
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 59963.43
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 59963.43
$ws.Range("N139").Value = -70243.42999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1086

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26267.5
$ws.Range("I54").Value = 15070
$ws.Range("J54").Value = 30000
$ws.Range("K54").Value = 15070
$ws.Range("L54").Value = 30000
$ws.Range("M54").Value = -14550
$ws.Range("N54").Value = -31040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12750.3
$ws.Range("I62").Value = 5200
$ws.Range("J62").Value = 17783.834
$ws.Range("K62").Value = 5200
$ws.Range("L62").Value = 17783.834
$ws.Range("M62").Value = -4576
$ws.Range("N62").Value = -19031.834

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 12750.3
$ws.Range("I65").Value = 5200
$ws.Range("J65").Value = 17783.834
$ws.Range("K65").Value = 26000
$ws.Range("L65").Value = 88919.17
$ws.Range("M65").Value = -22880
$ws.Range("N65").Value = -95159.17

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1577.8889
$ws.Range("I107").Value = 1944
$ws.Range("J107").Value = 296.5
$ws.Range("K107").Value = 5832
$ws.Range("L107").Value = 889.5
$ws.Range("M107").Value = -3912
$ws.Range("N107").Value = -4729.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 7001.3335
$ws.Range("I122").Value = 7001.3335
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 21004.0005
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -18554.0005
